$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Ref='D2'; Text='36.421.34'},
    @{Ref='E2'; Text='  -0.34%  '},
    @{Ref='D3'; Text='1.952.68'},
    @{Ref='E3'; Text='  -1.55%  '},
    @{Ref='E4'; Text='  -0.12%  '},
    @{Ref='D5'; Text='244.42'},
    @{Ref='E5'; Text='  -0.27%  '},
    @{Ref='D6'; Text='0.616'},
    @{Ref='E6'; Text='  -1.88%  '},
    @{Ref='B7'; Text='Solana'},
    @{Ref='C7'; Text='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'},
    @{Ref='D7'; Text='57.70'},
    @{Ref='E7'; Text='  -1.77%  '},
    @{Ref='B8'; Text='USDC'},
    @{Ref='C8'; Text='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'},
    @{Ref='D8'; Text='1.00'},
    @{Ref='E8'; Text='  -0.06%  '},
    @{Ref='D9'; Text='0.367'},
    @{Ref='E9'; Text='  -2.61%  '},
    @{Ref='D10'; Text='0.0850'},
    @{Ref='E10'; Text='  +4.34%  '},
    @{Ref='E11'; Text='  +0.38%  '},
    @{Ref='D12'; Text='2.239.67'},
    @{Ref='E12'; Text='  -1.55%  '},
    @{Ref='D13'; Text='0.823'},
    @{Ref='E13'; Text='  -5.22%  '},
    @{Ref='D14'; Text='21.47'},
    @{Ref='E14'; Text='  -11.32%  '},
    @{Ref='D15'; Text='13.57'},
    @{Ref='E15'; Text='  -4.05%  '},
    @{Ref='D16'; Text='5.20'},
    @{Ref='E16'; Text='  -4.80%  '},
    @{Ref='D17'; Text='1.950.57'},
    @{Ref='E17'; Text='  -1.88%  '},
    @{Ref='D18'; Text='36.355.68'},
    @{Ref='E18'; Text='  -0.21%  '},
    @{Ref='D19'; Text='0.0₃0885'},
    @{Ref='E19'; Text='  +2.18%  '},
    @{Ref='D20'; Text='69.83'},
    @{Ref='E20'; Text='  -1.92%  '},
    @{Ref='D21'; Text='230.12'},
    @{Ref='E21'; Text='  -2.14%  '},
    @{Ref='D22'; Text='5.08'},
    @{Ref='E22'; Text='  -4.82%  '},
    @{Ref='E23'; Text='  -0.01%  '},
    @{Ref='D24'; Text='2.43'},
    @{Ref='E24'; Text='  -7.31%  '},
    @{Ref='E25'; Text='  -0.55%  '},
    @{Ref='D26'; Text='9.30'},
    @{Ref='E26'; Text='  -9.41%  '},
    @{Ref='D27'; Text='161.69'},
    @{Ref='E28'; Text='  +6.77%  '},
    @{Ref='D29'; Text='19.43'},
    @{Ref='E29'; Text='  -2.34%  '},
    @{Ref='E30'; Text='  -1.77%  '},
    @{Ref='E31'; Text='  -2.02%  '},
    @{Ref='E32'; Text='  -5.36%  '},
    @{Ref='E33'; Text='  +3.26%  '},
    @{Ref='D34'; Text='4.29'},
    @{Ref='E34'; Text='  -4.08%  '},
    @{Ref='E35'; Text='  -0.66%  '},
    @{Ref='E36'; Text='  -0.18%  '},
    @{Ref='E37'; Text='  +1.06%  '},
    @{Ref='E38'; Text='  -5.48%  '},
    @{Ref='D39'; Text='3.04'},
    @{Ref='E39'; Text='  -1.63%  '},
    @{Ref='E40'; Text='  +1.41%  '},
    @{Ref='E41'; Text='  +0.43%  '},
    @{Ref='B42'; Text='VeChain'},
    @{Ref='C42'; Text='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'},
    @{Ref='D42'; Text='0.0212'},
    @{Ref='E42'; Text='  -1.11%  '},
    @{Ref='B43'; Text='TrustWalletToken'},
    @{Ref='C43'; Text='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'},
    @{Ref='D43'; Text='1.17'},
    @{Ref='E43'; Text='  -6.64%  '},
    @{Ref='D44'; Text='15.76'},
    @{Ref='E44'; Text='  -3.61%  '},
    @{Ref='E45'; Text='  -0.99%  '},
    @{Ref='E46'; Text='  -6.21%  '},
    @{Ref='D47'; Text='87.82'},
    @{Ref='E47'; Text='  -5.59%  '},
    @{Ref='D48'; Text='7.15'},
    @{Ref='E48'; Text='  -6.47%  '},
    @{Ref='E49'; Text='  -0.57%  '},
    @{Ref='D50'; Text='45.02'},
    @{Ref='E50'; Text='  -0.85%  '},
    @{Ref='D51'; Text='2.130.54'},
    @{Ref='E51'; Text='  -1.82%  '}
)

foreach ($change in $changes) {
    $cell = $ws.Range($change.Ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $change.Text
    $cell.Style = $origStyle
}
